# Logged Week 15 and simulated Week 16
# Append this period's per-play / per-game numbers onto the running
# season-long space-separated log strings, and bump the season totals
# on the summary tabs (OFF, DEF, ST, TURNS) accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS tab: per-play yardage logs (Home rush/pass, Road rush/pass)
# ---------------------------------------------------------------
$ydsSheet = $wb.Worksheets.Item("YDS")

$ydsSheet.Range("B2").Value = $ydsSheet.Range("B2").Value() + " 0 2 6 5 2 4 1 3 3 5 4 -1 -1 8 -6 2 30 0 1 0 13"
$ydsSheet.Range("C2").Value = $ydsSheet.Range("C2").Value() + " 4 5 4 4 4 6 5 8 1 4 19 0 9 4 4 10 3 1 8 8 38 1 9 12 7 -3 4 3 15 4 2 3 4 4 0 8 4 -3 1"
$ydsSheet.Range("B3").Value = $ydsSheet.Range("B3").Value() + " 19 7 6 46 5 1 7 5 9 7 5 12 29 5 6 1 13 10 3 -2"
$ydsSheet.Range("C3").Value = $ydsSheet.Range("C3").Value() + " 5 -1 19 28 16 45 6 5 8 17 21 20 5 9 7 9 7 34 19 17"

# ---------------------------------------------------------------
# OFF tab: season totals, Home (row 2) and Road (row 3)
# ---------------------------------------------------------------
$offSheet = $wb.Worksheets.Item("OFF")

$offSheet.Range("B2").Value = 10
$offSheet.Range("C2").Value = 350
$offSheet.Range("D2").Value = 23
$offSheet.Range("F2").Value = 140
$offSheet.Range("G2").Value = 102
$offSheet.Range("J2").Value = 72
$offSheet.Range("N2").Value = 40
$offSheet.Range("O2").Value = 47

$offSheet.Range("C3").Value = 386
$offSheet.Range("E3").Value = 60
$offSheet.Range("F3").Value = 196
$offSheet.Range("G3").Value = 63
$offSheet.Range("H3").Value = 52
$offSheet.Range("I3").Value = 126
$offSheet.Range("J3").Value = 106
$offSheet.Range("L3").Value = 517
$offSheet.Range("M3").Value = 342
$offSheet.Range("Q3").Value = 960

# ---------------------------------------------------------------
# DEF tab: season totals, Home (row 2) and Road (row 3)
# ---------------------------------------------------------------
$defSheet = $wb.Worksheets.Item("DEF")

$defSheet.Range("B2").Value = 7
$defSheet.Range("C2").Value = 355
$defSheet.Range("D2").Value = 22
$defSheet.Range("E2").Value = 16
$defSheet.Range("F2").Value = 122
$defSheet.Range("G2").Value = 102
$defSheet.Range("I2").Value = 9
$defSheet.Range("J2").Value = 59
$defSheet.Range("N2").Value = 32
$defSheet.Range("O2").Value = 37
$defSheet.Range("P2").Value = 16

$defSheet.Range("C3").Value = 350
$defSheet.Range("D3").Value = 5
$defSheet.Range("E3").Value = 64
$defSheet.Range("F3").Value = 201
$defSheet.Range("G3").Value = 70
$defSheet.Range("H3").Value = 46
$defSheet.Range("I3").Value = 102
$defSheet.Range("J3").Value = 133
$defSheet.Range("L3").Value = 521
$defSheet.Range("M3").Value = 346
$defSheet.Range("Q3").Value = 927

# ---------------------------------------------------------------
# ST tab: special-teams totals + per-kick logs
# ---------------------------------------------------------------
$stSheet = $wb.Worksheets.Item("ST")

$stSheet.Range("B2").Value = 140
$stSheet.Range("D2").Value = 124
$stSheet.Range("F2").Value = 96
$stSheet.Range("G2").Value = 84
$stSheet.Range("J2").Value = 45
$stSheet.Range("K2").Value = 43
$stSheet.Range("B3").Value = 107

$stSheet.Range("D3").Value = $stSheet.Range("D3").Value() + " 36 60 44 46 54"
$stSheet.Range("B4").Value = $stSheet.Range("B4").Value() + " 65 59 62"
$stSheet.Range("D4").Value = $stSheet.Range("D4").Value() + " 0 11 0 1 10"
$stSheet.Range("B5").Value = $stSheet.Range("B5").Value() + " 31 23 17"
$stSheet.Range("D5").Value = $stSheet.Range("D5").Value() + " 0 0"
$stSheet.Range("B6").Value = $stSheet.Range("B6").Value() + " 47 28"

# ---------------------------------------------------------------
# TURNS tab: Road (row 3) turnover totals
# ---------------------------------------------------------------
$turnsSheet = $wb.Worksheets.Item("TURNS")

$turnsSheet.Range("C3").Value = 8
$turnsSheet.Range("D3").Value = 16
$turnsSheet.Range("E3").Value = 18
